$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7568061351776123
$ws.Range("B1").Value = 1.959241390228271
$ws.Range("C1").Value = 4.057443618774414
$ws.Range("D1").Value = 3.562716484069824
$ws.Range("E1").Value = 2.007883310317993
